$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CaM binding to CaMKII is no longer a bimolecular reaction: ---
# give it (and the flanking flicker reaction) explicit rate-constant values.
$ws.Range("C6").Value = 20000
$ws.Range("C6").NumberFormat = "0.00E+00"

$ws.Range("F6").Value = 10000000
$ws.Range("F6").NumberFormat = "0.00E+00"

$ws.Range("C7").Value = 4200000
$ws.Range("C7").NumberFormat = "0.00E+00"

# Taller rows so the now-wrapped description text has room.
$ws.Rows("7").RowHeight = 29
$ws.Rows("8").RowHeight = 29

# --- Column A: narrower, left-aligned, word-wrapped text ---
$ws.Columns("A").ColumnWidth = 40.625

# Cells on the unshaded (white) rows.
$rWhite = $ws.Range("A1")
$rWhite = $excel.Union($rWhite, $ws.Range("A6"))
$rWhite = $excel.Union($rWhite, $ws.Range("A8"))
$rWhite = $excel.Union($rWhite, $ws.Range("A11"))
$rWhite.WrapText = $true
$rWhite.HorizontalAlignment = -4131

# Cells on the shaded rows (existing fill).
$rShade = $ws.Range("A2")
$rShade = $excel.Union($rShade, $ws.Range("A7"))
$rShade = $excel.Union($rShade, $ws.Range("A10"))
$rShade = $excel.Union($rShade, $ws.Range("A12"))
$rShade.WrapText = $true
$rShade.HorizontalAlignment = -4131

# New (empty) shaded cells continuing the band under row 2 - copy the
# finished A2 formatting so no extra intermediate styles get minted.
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selection moved by the author while editing ---
$ws.Range("G8").Select() | Out-Null
